# Auto-generated edit script: Add data for 2024-07-24
# Updates 2024 YTD (column K) violent-crime counts across Citywide Totals,
# By Neighborhood summary, and individual neighborhood sheets, plus a couple
# of small 2016 (column C) corrections that shipped in the same commit.

$wb = $excel.ActiveWorkbook

# Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 4431
$ws.Range("K3").Value = 4543
$ws.Range("C4").Value = 1064
$ws.Range("K4").Value = 914
$ws.Range("K6").Value = 5121
$ws.Range("C7").Value = 14980
$ws.Range("K7").Value = 15337

# By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 133
$ws.Range("K7").Value = 449
$ws.Range("K8").Value = 1024
$ws.Range("K11").Value = 301
$ws.Range("K14").Value = 86
$ws.Range("K15").Value = 155
$ws.Range("K19").Value = 465
$ws.Range("K20").Value = 351
$ws.Range("K29").Value = 808
$ws.Range("K30").Value = 55
$ws.Range("K31").Value = 167
$ws.Range("K33").Value = 641
$ws.Range("K34").Value = 79
$ws.Range("K36").Value = 193
$ws.Range("K37").Value = 521
$ws.Range("K42").Value = 568
$ws.Range("K43").Value = 138
$ws.Range("K44").Value = 136
$ws.Range("K47").Value = 96
$ws.Range("K48").Value = 199
$ws.Range("K49").Value = 88
$ws.Range("K51").Value = 195
$ws.Range("K52").Value = 406
$ws.Range("K53").Value = 204
$ws.Range("K54").Value = 287
$ws.Range("K60").Value = 99
$ws.Range("C63").Value = 156
$ws.Range("K63").Value = 48
$ws.Range("K64").Value = 94
$ws.Range("K65").Value = 347
$ws.Range("K67").Value = 592
$ws.Range("K68").Value = 39
$ws.Range("K72").Value = 70
$ws.Range("K75").Value = 53
$ws.Range("K76").Value = 212
$ws.Range("K79").Value = 387
$ws.Range("K83").Value = 328
$ws.Range("K85").Value = 688
$ws.Range("K88").Value = 178
$ws.Range("K89").Value = 221
$ws.Range("K90").Value = 142
$ws.Range("K93").Value = 57
$ws.Range("K95").Value = 273
$ws.Range("K96").Value = 171
$ws.Range("K98").Value = 77
$ws.Range("C101").Value = 14980
$ws.Range("K101").Value = 15337

# Bridgeport
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 86

# West Ridge
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 171

# Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 157
$ws.Range("K3").Value = 142
$ws.Range("K7").Value = 449

# Belmont Cragin
$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 99
$ws.Range("K7").Value = 301

# Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 221

# South Shore
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 241
$ws.Range("K3").Value = 230
$ws.Range("K7").Value = 688

# Little Village
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 153
$ws.Range("K7").Value = 406

# Logan Square
$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 204

# Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 307
$ws.Range("K4").Value = 59
$ws.Range("K6").Value = 345
$ws.Range("K7").Value = 1024

# South Chicago
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 113
$ws.Range("K3").Value = 119
$ws.Range("K7").Value = 328

# Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 238
$ws.Range("K6").Value = 186
$ws.Range("K7").Value = 641

# West Pullman
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 91
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 273

# Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 145
$ws.Range("K3").Value = 170
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 521

# New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 86
$ws.Range("K6").Value = 137
$ws.Range("K7").Value = 347

# Fuller Park
$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 55

# Gage Park
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 167

# North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 167
$ws.Range("K4").Value = 33
$ws.Range("K7").Value = 592

# Lincoln Park
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 88

# Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 287

# Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 287
$ws.Range("K6").Value = 226
$ws.Range("K7").Value = 808

# Lake View
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 199

# Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 144
$ws.Range("K3").Value = 142
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 465

# Irving Park
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 136

# River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 43
$ws.Range("K3").Value = 40
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 212

# Humboldt Park
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 153
$ws.Range("K3").Value = 175
$ws.Range("K4").Value = 23
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 568

# Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 127
$ws.Range("K7").Value = 387

# Near South Side
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 94

# Chicago Lawn
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 121
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 351

# Grand Boulevard
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 193

# West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 57

# Garfield Ridge
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 79

# Kenwood
$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 96

# Brighton Park
$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K4").Value = 11
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 155

# Wicker Park
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 77

# Albany Park
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 133

# United Center
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 178

# Pullman
$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 53

# Washington Heights
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 142

# Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 53
$ws.Range("K7").Value = 195

# North Park
$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 39

# Morgan Park
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 99

# Hyde Park
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 138

# Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 9
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 70
